$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 641.3333
$ws.Range("I2").Value = 524.75
$ws.Range("K2").Value = 524.75
$ws.Range("M2").Value = -411.75
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30228
$ws.Range("H28").Value = 126493.086
$ws.Range("I28").Value = 170094
$ws.Range("J28").Value = 534.8889
$ws.Range("K28").Value = 170094
$ws.Range("L28").Value = 534.8889
$ws.Range("M28").Value = -169609
$ws.Range("N28").Value = -1504.8889
$ws.Range("H40").Value = 10074.406
$ws.Range("I40").Value = 3038.1
$ws.Range("K40").Value = 3038.1
$ws.Range("M40").Value = -2863.1
$ws.Range("H62").Value = 669048.2
$ws.Range("I62").Value = 2001397
$ws.Range("K62").Value = 2001397
$ws.Range("M62").Value = -2000773
$ws.Range("H64").Value = 6452
$ws.Range("I64").Value = 4283.5
$ws.Range("J64").Value = 8310.714
$ws.Range("K64").Value = 4283.5
$ws.Range("L64").Value = 8310.714
$ws.Range("M64").Value = -4035.5
$ws.Range("N64").Value = -8806.714
$ws.Range("H65").Value = 669048.2
$ws.Range("I65").Value = 2001397
$ws.Range("K65").Value = 10006985
$ws.Range("M65").Value = -10003865
$ws.Range("H67").Value = 6452
$ws.Range("I67").Value = 4283.5
$ws.Range("J67").Value = 8310.714
$ws.Range("K67").Value = 4283.5
$ws.Range("L67").Value = 8310.714
$ws.Range("M67").Value = -3425.5
$ws.Range("N67").Value = -10026.714
$ws.Range("H70").Value = 113600.445
$ws.Range("I70").Value = 2700.8
$ws.Range("J70").Value = 252225
$ws.Range("K70").Value = 8102.400000000001
$ws.Range("L70").Value = 756675
$ws.Range("M70").Value = -7832.400000000001
$ws.Range("N70").Value = -757215
$ws.Range("H73").Value = 113600.445
$ws.Range("I73").Value = 2700.8
$ws.Range("J73").Value = 252225
$ws.Range("K73").Value = 8102.400000000001
$ws.Range("L73").Value = 756675
$ws.Range("M73").Value = -7166.400000000001
$ws.Range("N73").Value = -758547
$ws.Range("H98").Value = 973.05884
$ws.Range("I98").Value = 962.2857
$ws.Range("K98").Value = 962.2857
$ws.Range("M98").Value = 535.7143
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490
$ws.Range("H116").Value = 3443.6785
$ws.Range("I116").Value = 3346.476
$ws.Range("J116").Value = 3735.2856
$ws.Range("K116").Value = 3346.476
$ws.Range("L116").Value = 3735.2856
$ws.Range("M116").Value = 95.52399999999989
$ws.Range("N116").Value = -10619.2856
$ws.Range("H122").Value = 973.05884
$ws.Range("I122").Value = 962.2857
$ws.Range("K122").Value = 2886.8571
$ws.Range("M122").Value = -436.8571000000002
$ws.Range("H125").Value = 7782.1665
$ws.Range("I125").Value = 7369.5713
$ws.Range("J125").Value = 8359.799999999999
$ws.Range("K125").Value = 66326.14169999999
$ws.Range("L125").Value = 75238.2
$ws.Range("M125").Value = -63866.14169999999
$ws.Range("N125").Value = -80158.2
$ws.Range("H132").Value = 3123.275
$ws.Range("I132").Value = 3152.4167
$ws.Range("J132").Value = 2861
$ws.Range("K132").Value = 9457.250100000001
$ws.Range("L132").Value = 8583
$ws.Range("M132").Value = -6927.250100000001
$ws.Range("N132").Value = -13643
$ws.Range("H137").Value = 1681.5
$ws.Range("I137").Value = 1249.0625
$ws.Range("J137").Value = 1928.6072
$ws.Range("K137").Value = 3747.1875
$ws.Range("L137").Value = 5785.821599999999
$ws.Range("M137").Value = -1197.1875
$ws.Range("N137").Value = -10885.8216
$ws.Range("H138").Value = 6964.353
$ws.Range("J138").Value = 8270.421
$ws.Range("L138").Value = 24811.263
$ws.Range("N138").Value = -35091.263
$ws.Range("H141").Value = 2900.5908
$ws.Range("I141").Value = 3082.5625
$ws.Range("J141").Value = 2415.3333
$ws.Range("K141").Value = 9247.6875
$ws.Range("L141").Value = 7245.999899999999
$ws.Range("M141").Value = -4067.6875
$ws.Range("N141").Value = -17605.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2838.4707
$ws.Range("I2").Value = 2483.963
$ws.Range("K2").Value = 2483.963
$ws.Range("M2").Value = -2370.963
$ws.Range("H32").Value = 10142.692
$ws.Range("I32").Value = 952.3333
$ws.Range("K32").Value = 952.3333
$ws.Range("M32").Value = -665.3333
$ws.Range("H46").Value = 8261.25
$ws.Range("I46").Value = 9000
$ws.Range("J46").Value = 8015
$ws.Range("K46").Value = 9000
$ws.Range("L46").Value = 8015
$ws.Range("M46").Value = -8681
$ws.Range("N46").Value = -8653
$ws.Range("H61").Value = 4043.5806
$ws.Range("I61").Value = 3891.037
$ws.Range("K61").Value = 3891.037
$ws.Range("M61").Value = -3679.037
$ws.Range("H74").Value = 2438.5
$ws.Range("I74").Value = 2058
$ws.Range("K74").Value = 2058
$ws.Range("M74").Value = -1184
$ws.Range("H77").Value = 2438.5
$ws.Range("I77").Value = 2058
$ws.Range("K77").Value = 10290
$ws.Range("M77").Value = -5922
$ws.Range("H88").Value = 615.8333
$ws.Range("J88").Value = 491.875
$ws.Range("L88").Value = 491.875
$ws.Range("N88").Value = -1303.875
$ws.Range("H91").Value = 615.8333
$ws.Range("J91").Value = 491.875
$ws.Range("L91").Value = 491.875
$ws.Range("N91").Value = -3299.875
$ws.Range("H97").Value = 2344.625
$ws.Range("I97").Value = 813.44446
$ws.Range("K97").Value = 813.44446
$ws.Range("M97").Value = -317.44446
$ws.Range("H102").Value = 1780.875
$ws.Range("I102").Value = 1780.875
$ws.Range("K102").Value = 1780.875
$ws.Range("M102").Value = -158.875
$ws.Range("H110").Value = 939.3
$ws.Range("I110").Value = 954.96295
$ws.Range("J110").Value = 798.3333
$ws.Range("K110").Value = 954.96295
$ws.Range("L110").Value = 798.3333
$ws.Range("M110").Value = 1090.03705
$ws.Range("N110").Value = -4888.3333
$ws.Range("H116").Value = 2838.4707
$ws.Range("I116").Value = 2483.963
$ws.Range("K116").Value = 2483.963
$ws.Range("M116").Value = -189.9630000000002
$ws.Range("H122").Value = 1560.4348
$ws.Range("I122").Value = 1478.65
$ws.Range("K122").Value = 4435.950000000001
$ws.Range("M122").Value = -1985.950000000001
$ws.Range("H132").Value = 2656.7837
$ws.Range("I132").Value = 2633.6812
$ws.Range("J132").Value = 2975.6
$ws.Range("K132").Value = 7901.0436
$ws.Range("L132").Value = 8926.799999999999
$ws.Range("M132").Value = -5371.0436
$ws.Range("N132").Value = -13986.8
$ws.Range("H136").Value = 4043.5806
$ws.Range("I136").Value = 3891.037
$ws.Range("K136").Value = 11673.111
$ws.Range("M136").Value = -9123.110999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2838.4707
$ws.Range("I3").Value = 2483.963
$ws.Range("K3").Value = 2483.963
$ws.Range("M3").Value = -2369.963
$ws.Range("H5").Value = 2669.3333
$ws.Range("I5").Value = 2669.3333
$ws.Range("K5").Value = 2669.3333
$ws.Range("M5").Value = -2556.3333
$ws.Range("H19").Value = 4035
$ws.Range("I19").Value = 4035
$ws.Range("K19").Value = 4035
$ws.Range("M19").Value = -3862
$ws.Range("H86").Value = 5868
$ws.Range("I86").Value = 5798.5
$ws.Range("K86").Value = 5798.5
$ws.Range("M86").Value = -4675.5
$ws.Range("H89").Value = 5868
$ws.Range("I89").Value = 5798.5
$ws.Range("K89").Value = 28992.5
$ws.Range("M89").Value = -23376.5
$ws.Range("H107").Value = 1916.3784
$ws.Range("I107").Value = 1732.9286
$ws.Range("J107").Value = 2487.111
$ws.Range("K107").Value = 1732.9286
$ws.Range("L107").Value = 2487.111
$ws.Range("M107").Value = 187.0714
$ws.Range("N107").Value = -6327.111
$ws.Range("H134").Value = 1407.3334
$ws.Range("I134").Value = 1423.1578
$ws.Range("J134").Value = 1257
$ws.Range("K134").Value = 4269.4734
$ws.Range("L134").Value = 3771
$ws.Range("M134").Value = -1734.4734
$ws.Range("N134").Value = -8841

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 332.77777
$ws.Range("I7").Value = 314.2857
$ws.Range("J7").Value = 397.5
$ws.Range("K7").Value = 314.2857
$ws.Range("L7").Value = 397.5
$ws.Range("M7").Value = -201.2857
$ws.Range("N7").Value = -623.5
$ws.Range("H22").Value = 1742.6666
$ws.Range("I22").Value = 399.5
$ws.Range("J22").Value = 2011.3
$ws.Range("K22").Value = 399.5
$ws.Range("L22").Value = 2011.3
$ws.Range("M22").Value = -49.5
$ws.Range("N22").Value = -2711.3
$ws.Range("H31").Value = 5515.0625
$ws.Range("I31").Value = 8209.666999999999
$ws.Range("J31").Value = 3898.3
$ws.Range("K31").Value = 8209.666999999999
$ws.Range("L31").Value = 3898.3
$ws.Range("M31").Value = -7914.666999999999
$ws.Range("N31").Value = -4488.3
$ws.Range("H34").Value = 5515.0625
$ws.Range("I34").Value = 8209.666999999999
$ws.Range("J34").Value = 3898.3
$ws.Range("K34").Value = 8209.666999999999
$ws.Range("L34").Value = 3898.3
$ws.Range("M34").Value = -8007.666999999999
$ws.Range("N34").Value = -4302.3
$ws.Range("H58").Value = 8801.799999999999
$ws.Range("I58").Value = 9502.5
$ws.Range("J58").Value = 5999
$ws.Range("K58").Value = 9502.5
$ws.Range("L58").Value = 5999
$ws.Range("M58").Value = -9299.5
$ws.Range("N58").Value = -6405
$ws.Range("H62").Value = 6596.067
$ws.Range("I62").Value = 4897.375
$ws.Range("K62").Value = 4897.375
$ws.Range("M62").Value = -4273.375
$ws.Range("H65").Value = 6596.067
$ws.Range("I65").Value = 4897.375
$ws.Range("K65").Value = 24486.875
$ws.Range("M65").Value = -21366.875
$ws.Range("H132").Value = 5697.5713
$ws.Range("I132").Value = 5559.6665
$ws.Range("K132").Value = 16678.9995
$ws.Range("M132").Value = -14148.9995
$ws.Range("H134").Value = 5266.85
$ws.Range("I134").Value = 4372.8237
$ws.Range("J134").Value = 10333
$ws.Range("K134").Value = 13118.4711
$ws.Range("L134").Value = 30999
$ws.Range("M134").Value = -10583.4711
$ws.Range("N134").Value = -36069
$ws.Range("H136").Value = 8801.799999999999
$ws.Range("I136").Value = 9502.5
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 28507.5
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -25957.5
$ws.Range("N136").Value = -23097
$ws.Range("H137").Value = 69926.336
$ws.Range("J137").Value = 72499.5
$ws.Range("L137").Value = 72499.5
$ws.Range("N137").Value = -82699.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7401.353
$ws.Range("J55").Value = 7685.8125
$ws.Range("L55").Value = 23057.4375
$ws.Range("N55").Value = -23411.4375
$ws.Range("H56").Value = 9174.111000000001
$ws.Range("I56").Value = 9174.111000000001
$ws.Range("K56").Value = 9174.111000000001
$ws.Range("M56").Value = -8644.111000000001
$ws.Range("H68").Value = 2833.2856
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 3366.6
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 10099.8
$ws.Range("N68").Value = -11721.8
$ws.Range("M68").Value = -3689
$ws.Range("H71").Value = 2833.2856
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 3366.6
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 30299.4
$ws.Range("N71").Value = -38411.39999999999
$ws.Range("M71").Value = -9444
$ws.Range("H107").Value = 853.61536
$ws.Range("J107").Value = 687.375
$ws.Range("L107").Value = 2062.125
$ws.Range("N107").Value = -5902.125
$ws.Range("H131").Value = 2105.2144
$ws.Range("I131").Value = 979.8
$ws.Range("K131").Value = 2939.4
$ws.Range("M131").Value = 2100.6
$ws.Range("H132").Value = 3273.5
$ws.Range("I132").Value = 2631.3333
$ws.Range("J132").Value = 3915.6667
$ws.Range("K132").Value = 23681.9997
$ws.Range("L132").Value = 35241.0003
$ws.Range("M132").Value = -21151.9997
$ws.Range("N132").Value = -40301.0003
$ws.Range("H134").Value = 3782.6956
$ws.Range("I134").Value = 1361.15
$ws.Range("K134").Value = 4083.45
$ws.Range("M134").Value = 986.5499999999997
$ws.Range("H137").Value = 2950.7407
$ws.Range("I137").Value = 2517.0557
$ws.Range("J137").Value = 3818.111
$ws.Range("K137").Value = 7551.1671
$ws.Range("L137").Value = 11454.333
$ws.Range("M137").Value = -2451.1671
$ws.Range("N137").Value = -21654.333
$ws.Range("H138").Value = 8220.637000000001
$ws.Range("I138").Value = 3803.625
$ws.Range("K138").Value = 11410.875
$ws.Range("M138").Value = -6270.875
$ws.Range("H139").Value = 3522.0356
$ws.Range("I139").Value = 2234.8333
$ws.Range("J139").Value = 11245.25
$ws.Range("K139").Value = 6704.499899999999
$ws.Range("L139").Value = 33735.75
$ws.Range("M139").Value = -1564.499899999999
$ws.Range("N139").Value = -44015.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4522.222
$ws.Range("I70").Value = 4433.5
$ws.Range("K70").Value = 4433.5
$ws.Range("M70").Value = -4163.5
$ws.Range("H73").Value = 4522.222
$ws.Range("I73").Value = 4433.5
$ws.Range("K73").Value = 4433.5
$ws.Range("M73").Value = -3497.5
$ws.Range("H126").Value = 5057.364
$ws.Range("I126").Value = 4666.5
$ws.Range("K126").Value = 13999.5
$ws.Range("M126").Value = -11529.5
$ws.Range("H132").Value = 3099.6924
$ws.Range("I132").Value = 3183.72
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 9551.16
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -7021.16
$ws.Range("N132").Value = -8057
$ws.Range("H136").Value = 32709
$ws.Range("J136").Value = 32709
$ws.Range("L136").Value = 98127
$ws.Range("N136").Value = -103227

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17771.264
$ws.Range("I7").Value = 14070.272
$ws.Range("J7").Value = 22860.125
$ws.Range("K7").Value = 14070.272
$ws.Range("L7").Value = 22860.125
$ws.Range("M7").Value = -13958.272
$ws.Range("N7").Value = -23084.125
$ws.Range("H17").Value = 1806
$ws.Range("J17").Value = 1806
$ws.Range("L17").Value = 1806
$ws.Range("N17").Value = -2146
$ws.Range("H22").Value = 2238.5264
$ws.Range("I22").Value = 2193.1
$ws.Range("K22").Value = 2193.1
$ws.Range("M22").Value = -1898.1
$ws.Range("H27").Value = 2238.5264
$ws.Range("I27").Value = 2193.1
$ws.Range("K27").Value = 2193.1
$ws.Range("M27").Value = -2086.1
$ws.Range("H46").Value = 5882.1763
$ws.Range("I46").Value = 5736.636
$ws.Range("K46").Value = 5736.636
$ws.Range("M46").Value = -5548.636
$ws.Range("H61").Value = 3187.0476
$ws.Range("I61").Value = 3342.7334
$ws.Range("J61").Value = 2797.8333
$ws.Range("K61").Value = 3342.7334
$ws.Range("L61").Value = 2797.8333
$ws.Range("M61").Value = -3140.7334
$ws.Range("N61").Value = -3201.8333
$ws.Range("H68").Value = 2137.85
$ws.Range("I68").Value = 1964.0667
$ws.Range("J68").Value = 2659.2
$ws.Range("K68").Value = 1964.0667
$ws.Range("L68").Value = 2659.2
$ws.Range("M68").Value = -1215.0667
$ws.Range("N68").Value = -4157.2
$ws.Range("H71").Value = 2137.85
$ws.Range("I71").Value = 1964.0667
$ws.Range("J71").Value = 2659.2
$ws.Range("K71").Value = 9820.333500000001
$ws.Range("L71").Value = 13296
$ws.Range("M71").Value = -6076.333500000001
$ws.Range("N71").Value = -20784
$ws.Range("H113").Value = 3187.0476
$ws.Range("I113").Value = 3342.7334
$ws.Range("J113").Value = 2797.8333
$ws.Range("K113").Value = 3342.7334
$ws.Range("L113").Value = 2797.8333
$ws.Range("M113").Value = -1172.7334
$ws.Range("N113").Value = -7137.8333
$ws.Range("H126").Value = 17771.264
$ws.Range("I126").Value = 14070.272
$ws.Range("J126").Value = 22860.125
$ws.Range("K126").Value = 42210.81600000001
$ws.Range("L126").Value = 68580.375
$ws.Range("M126").Value = -39740.81600000001
$ws.Range("N126").Value = -73520.375
$ws.Range("H127").Value = 73444.164
$ws.Range("J127").Value = 73444.164
$ws.Range("L127").Value = 73444.164
$ws.Range("N127").Value = -83364.164
$ws.Range("H129").Value = 51500
$ws.Range("J129").Value = 51500
$ws.Range("L129").Value = 51500
$ws.Range("N129").Value = -61500
$ws.Range("H132").Value = 29309.393
$ws.Range("I132").Value = 35995.316
$ws.Range("K132").Value = 107985.948
$ws.Range("M132").Value = -105455.948
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 15000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20100
$ws.Range("H140").Value = 75264.39999999999
$ws.Range("J140").Value = 75264.39999999999
$ws.Range("L140").Value = 75264.39999999999
$ws.Range("N140").Value = -85624.39999999999
$ws.Range("H141").Value = 89996
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 89996
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 89996
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -100356

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 18936.875
$ws.Range("J69").Value = 20213.572
$ws.Range("L69").Value = 20213.572
$ws.Range("N69").Value = -21711.572
$ws.Range("H72").Value = 18936.875
$ws.Range("J72").Value = 20213.572
$ws.Range("L72").Value = 60640.716
$ws.Range("N72").Value = -68128.716
$ws.Range("H103").Value = 60000
$ws.Range("J103").Value = 60000
$ws.Range("L103").Value = 60000
$ws.Range("N103").Value = -62344
$ws.Range("H126").Value = 4082.6667
$ws.Range("I126").Value = 3399.2
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 10197.6
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -7727.599999999999
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 5091.875
$ws.Range("I132").Value = 4819.2856
$ws.Range("K132").Value = 14457.8568
$ws.Range("M132").Value = -11927.8568
$ws.Range("H136").Value = 12671.429
$ws.Range("I136").Value = 12671.429
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 38014.287
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("M136").Value = -35464.287
